# Natmi following Dr Hou advice
# Updates the LR-pairs (Cx3cl1-Itgb3) sheet with refreshed NATMI values and
# appends a new set of rows (sCs -> sCs sending/target block).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data block for rows 2-17 (columns A-T), matching the refreshed NATMI output.
$data = @(
    @('ECs', 'Cx3cl1', 'Itgb3', 'ECs', 3.0, 1.0, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3.0, 1.0, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 35.34803490008466, 318.132314100762, 0.1905765453263882, 0.1905765453263882),
    @('ECs', 'Cx3cl1', 'Itgb3', 'FAPs', 3.0, 1.0, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3.0, 1.0, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 40.53983464145111, 364.8585117730599, 0.2185677833551161, 0.2185677833551161),
    @('ECs', 'Cx3cl1', 'Itgb3', 'M2', 3.0, 1.0, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3.0, 1.0, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 2.722224632468444, 24.500021692216, 0.01467669044473503, 0.01467669044473503),
    @('ECs', 'Cx3cl1', 'Itgb3', 'sCs', 3.0, 1.0, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3.0, 1.0, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 4.503640927841777, 40.532768350576, 0.02428107621384519, 0.02428107621384519),
    @('FAPs', 'Cx3cl1', 'Itgb3', 'ECs', 3.0, 1.0, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3.0, 1.0, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 39.38559043034699, 354.470313873123, 0.2123447535647196, 0.2123447535647196),
    @('FAPs', 'Cx3cl1', 'Itgb3', 'FAPs', 3.0, 1.0, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3.0, 1.0, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 45.17041266410999, 406.5337139769899, 0.2435332323515613, 0.2435332323515613),
    @('FAPs', 'Cx3cl1', 'Itgb3', 'M2', 3.0, 1.0, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3.0, 1.0, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 3.033165060996, 27.298485548964, 0.01635310478682188, 0.01635310478682188),
    @('FAPs', 'Cx3cl1', 'Itgb3', 'sCs', 3.0, 1.0, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3.0, 1.0, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 5.018059915655999, 45.162539240904, 0.02705453148017167, 0.02705453148017167),
    @('sCs', 'Cx3cl1', 'Itgb3', 'ECs', 1.0, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3.0, 1.0, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 0.03733349627933334, 0.336001466514, 0.000201281026398836, 0.0002012810263988359),
    @('sCs', 'Cx3cl1', 'Itgb3', 'FAPs', 1.0, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3.0, 1.0, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 0.04281691386888889, 0.38535222482, 0.0002308445023813988, 0.0002308445023813988),
    @('sCs', 'Cx3cl1', 'Itgb3', 'M2', 1.0, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3.0, 1.0, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 0.002875129083555556, 0.025876161752, 0.00001550106447671664, 0.00001550106447671664),
    @('sCs', 'Cx3cl1', 'Itgb3', 'sCs', 1.0, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3.0, 1.0, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 0.004756605630222222, 0.042809450672, 0.00002564491833987715, 0.00002564491833987715),
    @('sCs', 'Cx3cl1', 'Itgb3', 'ECs', 3.0, 1.0, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3.0, 1.0, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 4.112927831434333, 37.016350482909, 0.02217457291493257, 0.02217457291493257),
    @('sCs', 'Cx3cl1', 'Itgb3', 'FAPs', 3.0, 1.0, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3.0, 1.0, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 4.717020752352221, 42.45318677116999, 0.02543149914152693, 0.02543149914152692),
    @('sCs', 'Cx3cl1', 'Itgb3', 'M2', 3.0, 1.0, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3.0, 1.0, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 0.3167450039568889, 2.850705035612, 0.001707709318894909, 0.001707709318894908),
    @('sCs', 'Cx3cl1', 'Itgb3', 'sCs', 3.0, 1.0, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3.0, 1.0, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 0.5240220614035555, 4.716198552631999, 0.002825229589689663, 0.002825229589689662),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
